$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 118-121 (shift data up by one match, add FT result columns) ---
# Row 118: match id 6814421
$ws.Range("A118").Value = 116
$ws.Range("B118").Value = 6814421
$ws.Range("C118").Value = "Slovenia Prva Liga"
$ws.Range("D118").Value = "Slovenia Prva Liga"
$ws.Range("E118").Value = 45360.45833333334
$ws.Range("F118").Value = "NK Bravo"
$ws.Range("G118").Value = "NK Aluminij"
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = "D"
$ws.Range("K118").Value = 1.666
$ws.Range("L118").Value = 3.5
$ws.Range("M118").Value = 5
$ws.Range("N118").Value = 1.5
$ws.Range("O118").Value = 3.8
$ws.Range("P118").Value = 6.5
$ws.Range("Q118").Value = -1
$ws.Range("R118").Value = 1.85
$ws.Range("S118").Value = 1.95
$ws.Range("T118").Value = 2.5
$ws.Range("U118").Value = 1.95
$ws.Range("V118").Value = 1.85
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = 2.8
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 0.95
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.8500000000000001

# Row 119: match id 6814423
$ws.Range("A119").Value = 117
$ws.Range("B119").Value = 6814423
$ws.Range("C119").Value = "Slovenia Prva Liga"
$ws.Range("D119").Value = "Slovenia Prva Liga"
$ws.Range("E119").Value = 45360.5625
$ws.Range("F119").Value = "NK Celje"
$ws.Range("G119").Value = "NK Maribor"
$ws.Range("H119").Value = 1
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = "D"
$ws.Range("K119").Value = 2
$ws.Range("L119").Value = 3.3
$ws.Range("M119").Value = 3.5
$ws.Range("N119").Value = 1.85
$ws.Range("O119").Value = 3.4
$ws.Range("P119").Value = 3.8
$ws.Range("Q119").Value = -0.5
$ws.Range("R119").Value = 1.85
$ws.Range("S119").Value = 1.95
$ws.Range("T119").Value = 2.5
$ws.Range("U119").Value = 1.85
$ws.Range("V119").Value = 1.95
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = 2.4
$ws.Range("Y119").Value = -1
$ws.Range("Z119").Value = -1
$ws.Range("AA119").Value = 0.95
$ws.Range("AB119").Value = -1
$ws.Range("AC119").Value = 0.95

# Row 120: match id 6816449
$ws.Range("A120").Value = 118
$ws.Range("B120").Value = 6816449
$ws.Range("C120").Value = "Slovenia Prva Liga"
$ws.Range("D120").Value = "Slovenia Prva Liga"
$ws.Range("E120").Value = 45361.375
$ws.Range("F120").Value = "NK Rogaska"
$ws.Range("G120").Value = "NK Radomlje"
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = "H"
$ws.Range("K120").Value = 2.625
$ws.Range("L120").Value = 3.2
$ws.Range("M120").Value = 2.5
$ws.Range("N120").Value = 2.45
$ws.Range("O120").Value = 3.2
$ws.Range("P120").Value = 2.7
$ws.Range("Q120").Value = 0
$ws.Range("R120").Value = 1.85
$ws.Range("S120").Value = 1.95
$ws.Range("T120").Value = 2.25
$ws.Range("U120").Value = 1.825
$ws.Range("V120").Value = 1.975
$ws.Range("W120").Value = 1.45
$ws.Range("X120").Value = -1
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.8500000000000001
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = 0.9750000000000001

# Row 121: match id 6814420
$ws.Range("A121").Value = 119
$ws.Range("B121").Value = 6814420
$ws.Range("C121").Value = "Slovenia Prva Liga"
$ws.Range("D121").Value = "Slovenia Prva Liga"
$ws.Range("E121").Value = 45361.45833333334
$ws.Range("F121").Value = "NS Mura"
$ws.Range("G121").Value = "Olimpija Ljubljana"
$ws.Range("H121").Value = 1
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = "D"
$ws.Range("K121").Value = 5.75
$ws.Range("L121").Value = 4
$ws.Range("M121").Value = 1.5
$ws.Range("N121").Value = 5.5
$ws.Range("O121").Value = 4
$ws.Range("P121").Value = 1.533
$ws.Range("Q121").Value = 1
$ws.Range("R121").Value = 1.95
$ws.Range("S121").Value = 1.85
$ws.Range("T121").Value = 2.75
$ws.Range("U121").Value = 1.85
$ws.Range("V121").Value = 1.95
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = 3
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = 0.95
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = -1
$ws.Range("AC121").Value = 0.95

# --- Add new rows 122-125 (new upcoming fixtures) ---
# Row 122: match id 6814427
$ws.Range("A121").Copy() | Out-Null
$ws.Range("A122").PasteSpecial(-4122) | Out-Null
$ws.Range("E121").Copy() | Out-Null
$ws.Range("E122").PasteSpecial(-4122) | Out-Null
$ws.Range("A122").Value = 120
$ws.Range("B122").Value = 6814427
$ws.Range("C122").Value = "Slovenia Prva Liga"
$ws.Range("D122").Value = "Slovenia Prva Liga"
$ws.Range("E122").Value = 45364.44791666666
$ws.Range("F122").Value = "NS Mura"
$ws.Range("G122").Value = "NK Bravo"
$ws.Range("K122").Value = 2.5
$ws.Range("L122").Value = 3.2
$ws.Range("M122").Value = 2.55
$ws.Range("N122").Value = 2.625
$ws.Range("O122").Value = 3.1
$ws.Range("P122").Value = 2.45
$ws.Range("Q122").Value = 0
$ws.Range("R122").Value = 1.975
$ws.Range("S122").Value = 1.825
$ws.Range("T122").Value = 2.25
$ws.Range("U122").Value = 1.975
$ws.Range("V122").Value = 1.825
$ws.Range("W122").Value = 0
$ws.Range("X122").Value = 0
$ws.Range("Y122").Value = 0
$ws.Range("Z122").Value = 0
$ws.Range("AA122").Value = 0

# Row 123: match id 6816448
$ws.Range("A121").Copy() | Out-Null
$ws.Range("A123").PasteSpecial(-4122) | Out-Null
$ws.Range("E121").Copy() | Out-Null
$ws.Range("E123").PasteSpecial(-4122) | Out-Null
$ws.Range("A123").Value = 121
$ws.Range("B123").Value = 6816448
$ws.Range("C123").Value = "Slovenia Prva Liga"
$ws.Range("D123").Value = "Slovenia Prva Liga"
$ws.Range("E123").Value = 45364.53125
$ws.Range("F123").Value = "NK Aluminij"
$ws.Range("G123").Value = "NK Rogaska"
$ws.Range("K123").Value = 2.15
$ws.Range("L123").Value = 3.2
$ws.Range("M123").Value = 3
$ws.Range("N123").Value = 2.6
$ws.Range("O123").Value = 3.2
$ws.Range("P123").Value = 2.4
$ws.Range("Q123").Value = 0
$ws.Range("R123").Value = 1.975
$ws.Range("S123").Value = 1.825
$ws.Range("T123").Value = 2.5
$ws.Range("U123").Value = 1.925
$ws.Range("V123").Value = 1.875
$ws.Range("W123").Value = 0
$ws.Range("X123").Value = 0
$ws.Range("Y123").Value = 0
$ws.Range("Z123").Value = 0
$ws.Range("AA123").Value = 0

# Row 124: match id 6814426
$ws.Range("A121").Copy() | Out-Null
$ws.Range("A124").PasteSpecial(-4122) | Out-Null
$ws.Range("E121").Copy() | Out-Null
$ws.Range("E124").PasteSpecial(-4122) | Out-Null
$ws.Range("A124").Value = 122
$ws.Range("B124").Value = 6814426
$ws.Range("C124").Value = "Slovenia Prva Liga"
$ws.Range("D124").Value = "Slovenia Prva Liga"
$ws.Range("E124").Value = 45364.61458333334
$ws.Range("F124").Value = "NK Radomlje"
$ws.Range("G124").Value = "NK Domzale"
$ws.Range("K124").Value = 2.55
$ws.Range("L124").Value = 3.25
$ws.Range("M124").Value = 2.4
$ws.Range("N124").Value = 2.25
$ws.Range("O124").Value = 3.3
$ws.Range("P124").Value = 2.7
$ws.Range("Q124").Value = -0.25
$ws.Range("R124").Value = 2.025
$ws.Range("S124").Value = 1.775
$ws.Range("T124").Value = 2.5
$ws.Range("U124").Value = 1.975
$ws.Range("V124").Value = 1.825
$ws.Range("W124").Value = 0
$ws.Range("X124").Value = 0
$ws.Range("Y124").Value = 0
$ws.Range("Z124").Value = 0
$ws.Range("AA124").Value = 0

# Row 125: match id 6814425
$ws.Range("A121").Copy() | Out-Null
$ws.Range("A125").PasteSpecial(-4122) | Out-Null
$ws.Range("E121").Copy() | Out-Null
$ws.Range("E125").PasteSpecial(-4122) | Out-Null
$ws.Range("A125").Value = 123
$ws.Range("B125").Value = 6814425
$ws.Range("C125").Value = "Slovenia Prva Liga"
$ws.Range("D125").Value = "Slovenia Prva Liga"
$ws.Range("E125").Value = 45365.53125
$ws.Range("F125").Value = "FC Koper"
$ws.Range("G125").Value = "NK Celje"
$ws.Range("K125").Value = 4
$ws.Range("L125").Value = 3.25
$ws.Range("M125").Value = 1.8
$ws.Range("N125").Value = 3.8
$ws.Range("O125").Value = 3.25
$ws.Range("P125").Value = 1.833
$ws.Range("Q125").Value = 0.5
$ws.Range("R125").Value = 1.9
$ws.Range("S125").Value = 1.9
$ws.Range("T125").Value = 2.5
$ws.Range("U125").Value = 1.95
$ws.Range("V125").Value = 1.85
$ws.Range("W125").Value = 0
$ws.Range("X125").Value = 0
$ws.Range("Y125").Value = 0
$ws.Range("Z125").Value = 0
$ws.Range("AA125").Value = 0

$excel.CutCopyMode = 0
